# Applies the recorded edits to the "Artfynd" worksheet:
#  - Row 62: corrected record id, sex, activity and public comment for the
#            Ceruchus chrysomelinus observation.
#  - Rows 63/65/66: the three fungi observations had been entered against
#    the wrong rows; rotate their data back into the correct rows
#    (new 63 <- old 66, new 65 <- old 63, new 66 <- old 65).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 62 - single-cell corrections
# ---------------------------------------------------------------------
$ws.Range("A62").Value = 80976096
$ws.Range("L62").Value = "hona"
$ws.Range("M62").Value = "funnen död"
$ws.Range("AC62").Value = "Fann ett huvud i perfekt rödmurken låga."

# ---------------------------------------------------------------------
# Rows 63 / 65 / 66 - rotate the mis-matched data back into place
# ---------------------------------------------------------------------

# Snapshot the "before" values of the three rows for every column that
# differs between them.
$cols = @("A","B","E","F","G","H","I","J","K","L","P","Q","R","S","Z","AB")

$row63 = @{}
$row65 = @{}
$row66 = @{}
foreach ($col in $cols) {
    $row63[$col] = $ws.Range($col + "63").Value2
    $row65[$col] = $ws.Range($col + "65").Value2
    $row66[$col] = $ws.Range($col + "66").Value2
}

# Column I ("Antal") is stored as text even when it looks numeric (e.g.
# "3"), so force text formatting on any destination cell about to receive
# a numeric-looking text value - otherwise Excel auto-coerces it to a
# number.
function Looks-Numeric($value) {
    if ($null -eq $value -or $value -eq "") { return $false }
    return ([string]$value) -match '^-?[0-9]+(\.[0-9]+)?$'
}

# new 63 <- old 66 ; new 65 <- old 63 ; new 66 <- old 65
# Skip any cell whose incoming (rotated) value is identical to what is
# already there - re-writing an already-empty cell with "" would turn a
# "blank text cell" into a truly-empty cell, which is a no-op for the
# source data but would show up as a spurious change.
foreach ($col in $cols) {
    if ($col -eq "I") {
        if (Looks-Numeric $row66["I"]) { $ws.Range("I63").NumberFormat = "@" }
        if (Looks-Numeric $row63["I"]) { $ws.Range("I65").NumberFormat = "@" }
        if (Looks-Numeric $row65["I"]) { $ws.Range("I66").NumberFormat = "@" }
    }
    if ($row63[$col] -ne $row66[$col]) { $ws.Range($col + "63").Value = $row66[$col] }
    if ($row65[$col] -ne $row63[$col]) { $ws.Range($col + "65").Value = $row63[$col] }
    if ($row66[$col] -ne $row65[$col]) { $ws.Range($col + "66").Value = $row65[$col] }
}
